# Auto-generated script to update cryptos price/volume columns
# Source: commit updating cryptos list on Sat Mar 30 08:45:09 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.748.68'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.22%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.499.79'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.04%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '597.42'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.71%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '194.17'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +4.82%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +1.36%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -3.03%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.649'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +1.30%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '53.65'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.06%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -2.59%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '9.49'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.32%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.055.81'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.98%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '606.44'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +5.59%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '69.926.05'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.05%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '18.94'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.32%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.54'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.70%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.495.80'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -1.03%  '
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.54%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.989'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.13%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '18.04'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +3.75%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '104.35'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +11.05%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +4.39%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.56'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -2.51%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.05'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +2.70%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.94'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.45%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.68'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +3.20%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '33.47'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +4.36%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.46'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +23.79%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.08'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +1.06%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '12.60'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +3.67%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +1.02%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '63.34'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.33%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.728.15'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +5.51%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0₃0805'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +3.46%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.04%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -7.13%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -2.86%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '36.56'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -1.31%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.98%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '500.78'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -6.41%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.135'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.32%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +1.00%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.32'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -2.36%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.41%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -4.16%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.42%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.71'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -4.14%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '131.84'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.000240'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.87%  '
